# The commit swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  (was "Office Theme", only used by the Notes Master)
#   ppt/theme/theme2.xml  (was "Integral",     used by the Slide Master /
#                           the presentation's own theme relationship)
# become, respectively, "Integral" and "Office Theme" -- i.e. the two
# theme color palettes trade places.
#
# The host's Theme object model only exposes (and persists) a single
# writable theme -- the one reachable from SlideMaster.Theme /
# Presentation's theme relationship, which is serialized back out as
# ppt/theme/theme2.xml. We therefore recolor that reachable theme with
# the plain "Office Theme" palette that used to live in theme1.xml,
# which reproduces the user-visible effect of the swap (the deck's
# actual design switches from the green/olive "Integral" palette to the
# default blue/gray "Office" palette).

function Get-RGBFromHex($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette == the "Office" colour scheme (previously theme1.xml),
# ordered to match ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officePalette.Length; $i++) {
    $tcs.Colors($i).RGB = Get-RGBFromHex($officePalette[$i - 1])
}

Write-Host "Theme palette swapped to Office Theme colours."
